# Finished Order testing automation
#
# 1. Orders sheet: mark every existing order (rows 2-15) as Delivered /
#    tested = TRUE in column E ("Test Result").
# 2. Update the on-screen selection/scroll state that Excel recorded when
#    the workbook was last saved: Customers scrolled down with G2 selected,
#    Orders scrolled down with E2 selected and made the active tab.

$wb = $excel.ActiveWorkbook

$wsCustomers = $wb.Worksheets.Item("Customers")
$wsOrders    = $wb.Worksheets.Item("Orders")

# --- Data change: populate the "Test Result" column on the Orders sheet ---
$wsOrders.Range("E2:E15").Value = $true

# --- View/selection state ---
# Customers: no longer the selected tab; scrolled so row 11 is at the top,
# with G2 selected.
$wsCustomers.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$wsCustomers.Range("G2").Select()

# Orders: becomes the active/selected tab, scrolled so row 10 is at the
# top, with E2 selected.
$wsOrders.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$wsOrders.Range("E2").Select()
